# Auto-generated script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Text)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (preventing numeric/locale auto-conversion of things like "25.38" or
    # "1.00"), then reset the cell style so no extra "quote prefix" number
    # format lingers on the cell.
    $Worksheet.Range($Address).Value2 = "'" + $Text
    $Worksheet.Range($Address).Style = "Normal"
}

Set-TextValue $ws "D2" "59.429.85"
Set-TextValue $ws "E2" "  -2.66%  "
Set-TextValue $ws "D3" "2.362.94"
Set-TextValue $ws "E3" "  -2.88%  "
Set-TextValue $ws "E4" "  +0.14%  "
Set-TextValue $ws "D5" "557.99"
Set-TextValue $ws "E5" "  -2.71%  "
Set-TextValue $ws "D6" "137.36"
Set-TextValue $ws "E6" "  -2.35%  "
Set-TextValue $ws "E7" "  -0.05%  "
Set-TextValue $ws "E8" "  -0.51%  "
Set-TextValue $ws "D9" "2.357.73"
Set-TextValue $ws "E9" "  -2.51%  "
Set-TextValue $ws "E10" "  -4.00%  "
Set-TextValue $ws "E11" "  -1.29%  "
Set-TextValue $ws "E13" "  -1.16%  "
Set-TextValue $ws "D14" "25.38"
Set-TextValue $ws "E14" "  -3.02%  "
Set-TextValue $ws "D15" "2.769.19"
Set-TextValue $ws "E15" "  -3.88%  "
Set-TextValue $ws "E16" "  -3.64%  "
Set-TextValue $ws "D17" "59.420.58"
Set-TextValue $ws "E17" "  -2.79%  "
Set-TextValue $ws "D18" "2.358.30"
Set-TextValue $ws "E18" "  -3.05%  "
Set-TextValue $ws "D19" "7.97"
Set-TextValue $ws "E19" "  +9.62%  "
Set-TextValue $ws "E20" "  -1.86%  "
Set-TextValue $ws "D21" "320.60"
Set-TextValue $ws "E21" "  -1.02%  "
Set-TextValue $ws "E22" "  -0.80%  "
Set-TextValue $ws "D23" "5.94"
Set-TextValue $ws "E23" "  -2.64%  "
Set-TextValue $ws "E24" "  +0.08%  "
Set-TextValue $ws "E25" "  -5.38%  "
Set-TextValue $ws "D26" "64.07"
Set-TextValue $ws "E26" "  -1.73%  "
Set-TextValue $ws "D27" "8.06"
Set-TextValue $ws "E27" "  -9.45%  "
Set-TextValue $ws "D28" "549.34"
Set-TextValue $ws "E28" "  -4.53%  "
Set-TextValue $ws "E29" "  -3.10%  "
Set-TextValue $ws "D30" "0.0₃0911"
Set-TextValue $ws "E30" "  -0.13%  "
Set-TextValue $ws "D31" "7.93"
Set-TextValue $ws "E31" "  +0.55%  "
Set-TextValue $ws "E32" "  -3.64%  "
Set-TextValue $ws "E33" "  -4.45%  "
Set-TextValue $ws "E34" "  -2.25%  "
Set-TextValue $ws "D35" "1.00"
Set-TextValue $ws "E35" "  -0.31%  "
Set-TextValue $ws "D36" "1.41"
Set-TextValue $ws "E36" "  +1.74%  "
Set-TextValue $ws "D37" "149.77"
Set-TextValue $ws "E37" "  -1.59%  "
Set-TextValue $ws "E38" "  -1.30%  "
Set-TextValue $ws "E39" "  -2.57%  "
Set-TextValue $ws "D40" "18.00"
Set-TextValue $ws "E40" "  -1.71%  "
Set-TextValue $ws "E41" "  -2.53%  "
Set-TextValue $ws "E42" "  +0.00%  "
Set-TextValue $ws "D43" "41.44"
Set-TextValue $ws "E43" "  -0.64%  "
Set-TextValue $ws "D44" "1.63"
Set-TextValue $ws "E44" "  -1.48%  "
Set-TextValue $ws "D45" "2.38"
Set-TextValue $ws "E45" "  +1.64%  "
Set-TextValue $ws "D46" "0.0₆0291"
Set-TextValue $ws "E46" "  +0.88%  "
Set-TextValue $ws "D47" "137.69"
Set-TextValue $ws "E47" "  -2.67%  "
Set-TextValue $ws "D48" "3.49"
Set-TextValue $ws "E48" "  -1.02%  "
Set-TextValue $ws "D49" "0.583"
Set-TextValue $ws "E49" "  -1.76%  "
Set-TextValue $ws "D50" "0.0498"
Set-TextValue $ws "E50" "  -1.94%  "
Set-TextValue $ws "D51" "18.95"
Set-TextValue $ws "E51" "  -3.06%  "
